# Adapt column header formatting to respective input file names:
#   "<col>_old" -> "<col>_FV2304"
#   "<col>_new" -> "<col>_FV2310"
# Wrap the header row + data range in an Excel Table ("Table1"), and
# freeze the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (A1:U1) -------------------------------------
$headers = @(
  "Segmentname_FV2304",
  "Segmentgruppe_FV2304",
  "Segment_FV2304",
  "Datenelement_FV2304",
  "Segment ID_FV2304",
  "Code_FV2304",
  "Qualifier_FV2304",
  "Beschreibung_FV2304",
  "Bedingungsausdruck_FV2304",
  "Bedingung_FV2304",
  "diff",
  "Segmentname_FV2310",
  "Segmentgruppe_FV2310",
  "Segment_FV2310",
  "Datenelement_FV2310",
  "Segment ID_FV2310",
  "Code_FV2310",
  "Qualifier_FV2310",
  "Beschreibung_FV2310",
  "Bedingungsausdruck_FV2310",
  "Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Turn the used range into a proper Excel Table ("Table1") ----------
$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $headers.Length
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"

# --- 3. Freeze the header row (row 1) --------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
